$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - bold/bordered header style matching existing headers (H1 has style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-20 for columns I and J
$values = @{
    2  = @(5, 7)
    3  = @(7, 8)
    4  = @(6, 8)
    5  = @(1, 3)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 5)
    11 = @(1, 7)
    12 = @(1, 4)
    13 = @(1, 7)
    14 = @(1, 7)
    15 = @(1, 6)
    16 = @(1, 5)
    17 = @(1, 4)
    18 = @(7, 8)
    19 = @(3, 5)
    20 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
